$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 6

$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 10

$ws.Activate()
$ws.Range("D5").Select()
